# Apply NATMI recalculated values to Ptprz1-L1cam sheet (per Dr Hou's advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1030763333333333
$ws.Range("H2").Value = 0.309229
$ws.Range("I2").Value = 0.01126512502660735
$ws.Range("J2").Value = 0.01126512502660735
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.59487733333333
$ws.Range("N2").Value = 67.784632
$ws.Range("O2").Value = 0.7395019553569895
$ws.Range("P2").Value = 0.7395019553569895
$ws.Range("Q2").Value = 2.328997107636444
$ws.Range("R2").Value = 20.960973968728
$ws.Range("S2").Value = 0.008330581984517095
$ws.Range("T2").Value = 0.008330581984517096

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1030763333333333
$ws.Range("H3").Value = 0.309229
$ws.Range("I3").Value = 0.01126512502660735
$ws.Range("J3").Value = 0.01126512502660735
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.3045986666666667
$ws.Range("N3").Value = 0.913796
$ws.Range("O3").Value = 0.00996913177602551
$ws.Range("P3").Value = 0.00996913177602551
$ws.Range("Q3").Value = 0.03139691369822222
$ws.Range("R3").Value = 0.282572223284
$ws.Range("S3").Value = 0.0001123035158636516
$ws.Range("T3").Value = 0.0001123035158636516

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1030763333333333
$ws.Range("H4").Value = 0.309229
$ws.Range("I4").Value = 0.01126512502660735
$ws.Range("J4").Value = 0.01126512502660735
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.654706
$ws.Range("N4").Value = 22.964118
$ws.Range("O4").Value = 0.2505289128669849
$ws.Range("P4").Value = 0.2505289128669849
$ws.Range("Q4").Value = 0.7890190272246667
$ws.Range("R4").Value = 7.101171245021999
$ws.Range("S4").Value = 0.002822239526226605
$ws.Range("T4").Value = 0.002822239526226605

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.046962666666667
$ws.Range("H5").Value = 27.140888
$ws.Range("I5").Value = 0.9887348749733926
$ws.Range("J5").Value = 0.9887348749733927
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.59487733333333
$ws.Range("N5").Value = 67.784632
$ws.Range("O5").Value = 0.7395019553569895
$ws.Range("P5").Value = 0.7395019553569895
$ws.Range("Q5").Value = 204.4150116925796
$ws.Range("R5").Value = 1839.735105233216
$ws.Range("S5").Value = 0.7311713733724724
$ws.Range("T5").Value = 0.7311713733724725

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.046962666666667
$ws.Range("H6").Value = 27.140888
$ws.Range("I6").Value = 0.9887348749733926
$ws.Range("J6").Value = 0.9887348749733927
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.3045986666666667
$ws.Range("N6").Value = 0.913796
$ws.Range("O6").Value = 0.00996913177602551
$ws.Range("P6").Value = 0.00996913177602551
$ws.Range("Q6").Value = 2.755692765649778
$ws.Range("R6").Value = 24.801234890848
$ws.Range("S6").Value = 0.009856828260161859
$ws.Range("T6").Value = 0.009856828260161859

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.046962666666667
$ws.Range("H7").Value = 27.140888
$ws.Range("I7").Value = 0.9887348749733926
$ws.Range("J7").Value = 0.9887348749733927
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.654706
$ws.Range("N7").Value = 22.964118
$ws.Range("O7").Value = 0.2505289128669849
$ws.Range("P7").Value = 0.2505289128669849
$ws.Range("Q7").Value = 69.25183940630934
$ws.Range("R7").Value = 623.266554656784
$ws.Range("S7").Value = 0.2477066733407583
$ws.Range("T7").Value = 0.2477066733407583

